$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16, shifting existing rows 16-20 down to 17-21
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44438
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100108
$ws.Cells.Item(16, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(16, 9).Value = 100108003
$ws.Cells.Item(16, 10).Value = "Maracuyá"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 25
$ws.Cells.Item(16, 14).Value = 35000
$ws.Cells.Item(16, 15).Value = 35000
$ws.Cells.Item(16, 16).Value = 35000
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(16, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 19).Value = 1944
$ws.Cells.Item(16, 20).Value = 18
